# Update the "today" date placeholder that is cached on every slide layout
# (Insert > Header & Footer > Date and time, fixed text) from 02/04/2025 to
# 06/04/2025. The placeholder is a field, but it is re-typed as a normal
# date string here (same as re-entering the text by hand), so we simply
# look up the date placeholder on every layout of the slide master and
# rewrite its text.
$p = $ppt.ActivePresentation
$m = $p.SlideMaster

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $cl = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "06/04/2025"
        }
    }
}

# Nudge the three signature rectangles on the single content slide a touch
# to the right (layout tweak), keeping their vertical position and size.
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    switch ($sh.Name) {
        "Retângulo 4"  { $sh.Left = 389.9881286621094 }   # 4935597 -> 4952849 EMU
        "Retângulo 16" { $sh.Left = 568.1328735351562 }   # 7198035 -> 7215287 EMU
        "Retângulo 10" { $sh.Left = 211.49551391601562 }  # 2668741 -> 2685993 EMU
    }
}
